# Update countries & provincias Spain
#
# 1) Swap the "Croacia"/"Grecia" row contents (row 91 <-> row 92, country
#    name column A only) plus refresh their statistics.
# 2) Refresh COVID-19 statistics for a handful of other countries.
# 3) Bump the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Croacia / Grecia labels (A91 <-> A92) -----------------------
$ws.Range("A91").Value = "Grecia"
$ws.Range("A92").Value = "Croacia"

# --- Refresh numeric statistics (B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy,
#     H=Muertes) ----------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6218003
$ws.Range("C4").Value = 6207
$ws.Range("D4").Value = 3458559
$ws.Range("E4").Value = 2571542
$ws.Range("G4").Value = 166
$ws.Range("H4").Value = 187902

# Row 6 - India
$ws.Range("B6").Value = 3733936
$ws.Range("C6").Value = 45997
$ws.Range("D6").Value = 2872714
$ws.Range("E6").Value = 795300
$ws.Range("G6").Value = 487
$ws.Range("H6").Value = 65922

# Row 14 - Chile
$ws.Range("B14").Value = 413145
$ws.Range("C14").Value = 1419
$ws.Range("D14").Value = 385790
$ws.Range("E14").Value = 16034
$ws.Range("G14").Value = 32
$ws.Range("H14").Value = 11321

# Row 16 - Reino Unido
$ws.Range("B16").Value = 337168
$ws.Range("C16").Value = 1295
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 41504

# Row 23 - Alemania
$ws.Range("B23").Value = 245408
$ws.Range("C23").Value = 616
$ws.Range("E23").Value = 16131
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9377

# Row 46 - Emiratos Arabes Unidos
$ws.Range("B46").Value = 70805
$ws.Range("C46").Value = 574
$ws.Range("D46").Value = 61491
$ws.Range("E46").Value = 8930

# Row 52 - Singapur
$ws.Range("D52").Value = 55749
$ws.Range("E52").Value = 1076

# Row 66 - Moldavia
$ws.Range("B66").Value = 37208
$ws.Range("C66").Value = 288
$ws.Range("E66").Value = 10328

# Row 90 - Noruega
$ws.Range("B90").Value = 10840
$ws.Range("C90").Value = 58
$ws.Range("E90").Value = 1228

# Row 91 - now Grecia (was Croacia)
$ws.Range("B91").Value = 10524
$ws.Range("C91").Value = 207
$ws.Range("D91").Value = 3804
$ws.Range("E91").Value = 6449
$ws.Range("G91").Value = 5
$ws.Range("H91").Value = 271

# Row 92 - now Croacia (was Grecia)
$ws.Range("B92").Value = 10414
$ws.Range("C92").Value = 145
$ws.Range("D92").Value = 7735
$ws.Range("E92").Value = 2492
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 187

# Row 94 - Albania
$ws.Range("B94").Value = 9606
$ws.Range("C94").Value = 93
$ws.Range("D94").Value = 5441
$ws.Range("E94").Value = 3875
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = 290

# Row 103 - Namibia
$ws.Range("B103").Value = 7692
$ws.Range("C103").Value = 142
$ws.Range("E103").Value = 4284
$ws.Range("G103").Value = 6
$ws.Range("H103").Value = 81

# Row 129 - Gambia
$ws.Range("B129").Value = 3029
$ws.Range("C129").Value = 66
$ws.Range("E129").Value = 1901

# Row 147 - Trinidad yTobago
$ws.Range("B147").Value = 1773
$ws.Range("C147").Value = 14
$ws.Range("D147").Value = 685
$ws.Range("E147").Value = 1063
$ws.Range("G147").Value = 3
$ws.Range("H147").Value = 25

# Row 150 - Reunion
$ws.Range("B150").Value = 1714
$ws.Range("C150").Value = 35
$ws.Range("E150").Value = 825

# --- Update the "last updated" timestamp string ------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 17:49"
